$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Register "Date 14-Mar" into the shared string table first (it ends up
# earlier in the table than "Date 14-Mar-12"), by writing it to A17 first.
$ws.Range("A17").Value = "Date 14-Mar"
$ws.Range("A16").Value = "Date 14-Mar-12"

# Apply the same label styling (grey font, same as other "Date ..." rows A7:A15)
$ws.Range("A7").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# B17: "Date 14-Mar" value, formatted with built-in d-mmm (numFmt 16)
$ws.Range("B17").Value = 40982.563138888887
$ws.Range("B17").NumberFormat = "d-mmm"

# B16: "Date 14-Mar-12" value, formatted with built-in d-mmm-yy (numFmt 15)
$ws.Range("B16").Value = 40982.563138888887
$ws.Range("B16").NumberFormat = "d-mmm-yy"

# Extend the shared formula in column C down through the new rows
$ws.Range("C16").Formula = "=B16"
$ws.Range("C16").NumberFormat = $ws.Range("C15").NumberFormat

$ws.Range("C17").Formula = "=B17"
$ws.Range("C17").NumberFormat = $ws.Range("C15").NumberFormat

$ws.Range("A17").Select()
